# The post about "「オレンジの効能が何か知っていますか？」" (row 667) was removed.
# Deleting the entire row shifts all subsequent rows (668-699) up by one,
# which also updates the sheet dimension from A1:C699 to A1:C698 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("667").Delete()
